$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, [string]$Text)
    if ($Text -match '^-?\d+(\.\d+)?$') {
        $Cell.Value = "'" + $Text
        $Cell.Style = "Normal"
    } else {
        $Cell.Value = $Text
    }
}


Set-TextValue $ws.Range("D2") '69.285.53'
Set-TextValue $ws.Range("E2") '  -0.13%  '
Set-TextValue $ws.Range("D3") '3.672.63'
Set-TextValue $ws.Range("E3") '  -0.40%  '
Set-TextValue $ws.Range("E4") '  +0.19%  '
Set-TextValue $ws.Range("D5") '681.52'
Set-TextValue $ws.Range("E5") '  -0.26%  '
Set-TextValue $ws.Range("D6") '157.88'
Set-TextValue $ws.Range("E6") '  -3.04%  '
Set-TextValue $ws.Range("E7") '  +0.06%  '
Set-TextValue $ws.Range("E8") '  -1.36%  '
Set-TextValue $ws.Range("D9") '0.145'
Set-TextValue $ws.Range("E9") '  -1.90%  '
Set-TextValue $ws.Range("E10") '  -3.74%  '
Set-TextValue $ws.Range("D11") '0.435'
Set-TextValue $ws.Range("E11") '  -3.15%  '
Set-TextValue $ws.Range("E12") '  -2.53%  '
Set-TextValue $ws.Range("D13") '4.295.85'
Set-TextValue $ws.Range("E13") '  -0.31%  '
Set-TextValue $ws.Range("D14") '32.14'
Set-TextValue $ws.Range("E14") '  -4.39%  '
Set-TextValue $ws.Range("D15") '3.676.68'
Set-TextValue $ws.Range("E15") '  -0.32%  '
Set-TextValue $ws.Range("D16") '69.256.76'
Set-TextValue $ws.Range("E16") '  -0.20%  '
Set-TextValue $ws.Range("E17") '  +1.79%  '
Set-TextValue $ws.Range("D18") '15.88'
Set-TextValue $ws.Range("E19") '  -4.14%  '
Set-TextValue $ws.Range("D20") '468.99'
Set-TextValue $ws.Range("E20") '  -2.27%  '
Set-TextValue $ws.Range("D21") '9.94'
Set-TextValue $ws.Range("E21") '  +1.22%  '
Set-TextValue $ws.Range("D22") '0.647'
Set-TextValue $ws.Range("E22") '  -3.05%  '
Set-TextValue $ws.Range("D23") '79.96'
Set-TextValue $ws.Range("E23") '  -0.21%  '
Set-TextValue $ws.Range("D24") '3.819.85'
Set-TextValue $ws.Range("E24") '  -0.38%  '
Set-TextValue $ws.Range("E25") '  -0.05%  '
Set-TextValue $ws.Range("D26") '0.0000120'
Set-TextValue $ws.Range("E26") '  -6.01%  '
Set-TextValue $ws.Range("D27") '10.88'
Set-TextValue $ws.Range("E27") '  -5.26%  '
Set-TextValue $ws.Range("D28") '9.09'
Set-TextValue $ws.Range("E28") '  -5.23%  '
Set-TextValue $ws.Range("D29") '2.69'
Set-TextValue $ws.Range("E29") '  -2.10%  '
Set-TextValue $ws.Range("D30") '1.73'
Set-TextValue $ws.Range("E30") '  -5.30%  '
Set-TextValue $ws.Range("D33") '1.98'
Set-TextValue $ws.Range("E33") '  -6.29%  '
Set-TextValue $ws.Range("D34") '26.79'
Set-TextValue $ws.Range("E34") '  -1.20%  '
Set-TextValue $ws.Range("D35") '3.653.01'
Set-TextValue $ws.Range("E35") '  +0.06%  '
Set-TextValue $ws.Range("E36") '  -3.85%  '
Set-TextValue $ws.Range("D37") '8.17'
Set-TextValue $ws.Range("E37") '  -4.82%  '
Set-TextValue $ws.Range("D38") '6.02'
Set-TextValue $ws.Range("E38") '  -2.90%  '
Set-TextValue $ws.Range("E39") '  -0.01%  '
Set-TextValue $ws.Range("D40") '2.23'
Set-TextValue $ws.Range("E40") '  +2.85%  '
Set-TextValue $ws.Range("D41") '0.0897'
Set-TextValue $ws.Range("E41") '  -4.86%  '
Set-TextValue $ws.Range("E42") '  +0.03%  '
Set-TextValue $ws.Range("D43") '167.19'
Set-TextValue $ws.Range("E43") '  +7.93%  '
Set-TextValue $ws.Range("D44") '0.939'
Set-TextValue $ws.Range("E44") '  -2.22%  '
Set-TextValue $ws.Range("D45") '47.58'
Set-TextValue $ws.Range("E45") '  -1.35%  '
Set-TextValue $ws.Range("D46") '2.70'
Set-TextValue $ws.Range("E46") '  -5.58%  '
Set-TextValue $ws.Range("D47") '0.000275'
Set-TextValue $ws.Range("E47") '  -2.50%  '
Set-TextValue $ws.Range("E48") '  +1.43%  '
Set-TextValue $ws.Range("D49") '1.26'
Set-TextValue $ws.Range("E49") '  -4.79%  '
Set-TextValue $ws.Range("D50") '7.73'
Set-TextValue $ws.Range("E50") '  -4.56%  '
Set-TextValue $ws.Range("D51") '26.78'
Set-TextValue $ws.Range("E51") '  -3.51%  '

# Row 31 / Row 32: NEARProtocol and Binance-PegBSC-USD swap rows, with updated values
Set-TextValue $ws.Range("B31") 'Binance-PegBSC-USD'
Set-TextValue $ws.Range("C31") 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue $ws.Range("D31") '1.00'
Set-TextValue $ws.Range("E31") '  +0.01%  '

Set-TextValue $ws.Range("B32") 'NEARProtocol'
Set-TextValue $ws.Range("C32") 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Range("D32") '6.54'
Set-TextValue $ws.Range("E32") '  -4.45%  '
